# Daily attendance processing - 2026-01-11 14:57:36
# Swap the order of "dnasr281@gmail.com" and "System" in the "Recorded By"
# column (G) for every row where both names are present together.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
